$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case2")

# Fix the expected-response text: "用户名或者密码错误" -> "用户名或密码错误"
# (drop the extra "者") for both affected test cases (rows 10 and 11).
$ws.Range("E10").Value = '{"status":0,"code":"20111","data":null,"msg":"用户名或密码错误"}'
$ws.Range("E11").Value = '{"status":0,"code":"20111","data":null,"msg":"用户名或密码错误"}'

# Reflect the reviewer's final selection state on the Case2 sheet.
$ws.Range("F2:F12").Select() | Out-Null
